# Update RAD Test Case data on Sheet1 of TaxPayerSSNNoMatch.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Refresh timestamps for the four test rows (new test run results)
$ws.Range("B2").Value = "Fri Sep 08 18:10:19 EDT 2023"
$ws.Range("B3").Value = "Fri Sep 08 18:10:34 EDT 2023"
$ws.Range("B4").Value = "Fri Sep 08 18:10:48 EDT 2023"
$ws.Range("B5").Value = "Fri Sep 08 18:11:03 EDT 2023"

# Fix payment type label (singular -> plural)
$ws.Range("D4").Value = "Extension Payments"

# Update the active selection to D4
$ws.Range("D4").Select()
